$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Harry Potter and the Goblet of Fire ---
$ws.Range("C11:D11").Copy($ws.Range("C12:D12"))
$ws.Range("A12").Value = "Harry Potter and the Goblet of Fire"
$ws.Range("B12").Value = "J.K. Rowling"
$ws.Range("C12").Value = 44217
$ws.Range("D12").Value = 44224
$ws.Range("E12").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F12").Value = "Audio"
$ws.Range("G12").Value = "21 Hours 29 Mins"
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = $true

# --- Row 13: Harry Potter and the Order of the Phoenix ---
$ws.Range("C11:D11").Copy($ws.Range("C13:D13"))
$ws.Range("C11").Copy($ws.Range("B13"))
$ws.Range("C11").Copy($ws.Range("F13"))
$ws.Range("C11").Copy($ws.Range("G13"))
$ws.Range("A13").Value = "Harry Potter and the Order of the Phoenix"
$ws.Range("B13").Value = "J.K. Rowling"
$ws.Range("C13").Value = 44224
$ws.Range("D13").Value = 44232
$ws.Range("E13").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F13").Value = "Audio"
$ws.Range("G13").Value = "27 Hours 19 Mins"
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = $true

# --- Row 14: Harry Potter and the Half Blood Prince ---
$ws.Range("C11:D11").Copy($ws.Range("C14:D14"))
$ws.Range("A14").Value = "Harry Potter and the Half Blood Prince"
$ws.Range("B14").Value = "J.K. Rowling"
$ws.Range("C14").Value = 44232
$ws.Range("D14").Value = 44240
$ws.Range("E14").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F14").Value = "Audio"
$ws.Range("G14").Value = "19 Hours 6 Mins"
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = $true

# --- Row 15: Wizard's First Rule ---
$ws.Range("C11:D11").Copy($ws.Range("C15:D15"))
$ws.Range("A15").Value = "Wizard's First Rule"
$ws.Range("B15").Value = "Terry Goodkind"
$ws.Range("C15").Value = 44214
$ws.Range("D15").Value = 44240
$ws.Range("E15").Value = "fiction;wizards;evil;good;seeker"
$ws.Range("F15").Value = "Hard Copy"
$ws.Range("G15").Value = "820 Pages"
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = $true

$ws.Range("A16").Select()

# Column A is widened (best-fit) to accommodate the longer book titles just added.
$ws.Columns("A").ColumnWidth = 34.6
